# The "第三组项目计划表" workbook tracks each team member's task status in
# column C ("完成情况"). This upload fills in the "完成" (Done) status for
# the first week's five tasks (rows 3-8), widens the plan-content column
# (B) so the longer task descriptions are fully visible, and leaves the
# cursor on the newly-entered C4 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the first week's tasks as completed.
$ws.Range("C5").Value = "完成"
$ws.Range("C6").Value = "完成"
$ws.Range("C7").Value = "完成"
$ws.Range("C8").Value = "完成"
$ws.Range("C3").Value = "完成"
$ws.Range("C4").Value = "完成"

# Widen column B ("计划内容") to better fit the task descriptions.
# (target stored width 54.125 chars; Excel quantizes ColumnWidth to whole
# pixels at the default Normal-style digit width, so 374/7 rounds back to
# the stored width nearest 54.125 after that pixel snap.)
$ws.Columns.Item(2).ColumnWidth = 53.4285714285714

# Leave the selection on C4, matching where the edit finished.
[void]$ws.Range("C4").Select()
